$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "keyword" column (B) for rows 36-38 from "bitcoin" to "best bitcoin"
$ws.Range("B36").Value = "best bitcoin"
$ws.Range("B37").Value = "best bitcoin"
$ws.Range("B38").Value = "best bitcoin"

# Move the active selection to B36
$ws.Range("B36").Select()
